$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1305.2307
$ws.Range("I53").Value = 688.8333
$ws.Range("J53").Value = 1833.5714
$ws.Range("K53").Value = 688.8333
$ws.Range("L53").Value = 1833.5714
$ws.Range("M53").Value = -51.83330000000001
$ws.Range("N53").Value = -3107.5714
$ws.Range("H100").Value = 4578.3
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 4722.875
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 4722.875
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -5804.875
$ws.Range("H125").Value = 4939.8
$ws.Range("J125").Value = 5925
$ws.Range("L125").Value = 53325
$ws.Range("N125").Value = -58245
$ws.Range("H132").Value = 1383.3
$ws.Range("I132").Value = 1379.8853
$ws.Range("K132").Value = 4139.6559
$ws.Range("M132").Value = -1609.6559
$ws.Range("H137").Value = 1853773.2
$ws.Range("I137").Value = 1652.7778
$ws.Range("J137").Value = 4631954
$ws.Range("K137").Value = 4958.3334
$ws.Range("L137").Value = 13895862
$ws.Range("M137").Value = -2408.3334
$ws.Range("N137").Value = -13900962
$ws.Range("H138").Value = 2542.36
$ws.Range("I138").Value = 1296.0769
$ws.Range("J138").Value = 2728.5862
$ws.Range("K138").Value = 3888.2307
$ws.Range("L138").Value = 8185.758600000001
$ws.Range("M138").Value = 1251.7693
$ws.Range("N138").Value = -18465.7586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 4990
$ws.Range("I41").Value = 4990
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4990
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4576
$ws.Range("H61").Value = 3357.5
$ws.Range("I61").Value = 3119.375
$ws.Range("K61").Value = 3119.375
$ws.Range("M61").Value = -2907.375
$ws.Range("H74").Value = 2750.6667
$ws.Range("I74").Value = 2456.2964
$ws.Range("K74").Value = 2456.2964
$ws.Range("M74").Value = -1582.2964
$ws.Range("H77").Value = 2750.6667
$ws.Range("I77").Value = 2456.2964
$ws.Range("K77").Value = 12281.482
$ws.Range("M77").Value = -7913.482
$ws.Range("H110").Value = 2267.75
$ws.Range("I110").Value = 2179.1538
$ws.Range("K110").Value = 2179.1538
$ws.Range("M110").Value = -134.1538
$ws.Range("H132").Value = 191701.64
$ws.Range("I132").Value = 240254.38
$ws.Range("K132").Value = 720763.14
$ws.Range("M132").Value = -718233.14
$ws.Range("H136").Value = 3357.5
$ws.Range("I136").Value = 3119.375
$ws.Range("K136").Value = 9358.125
$ws.Range("M136").Value = -6808.125
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 69989
$ws.Range("J52").Value = 69989
$ws.Range("L52").Value = 69989
$ws.Range("N52").Value = -70515
$ws.Range("H99").Value = 3499.889
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H107").Value = 1426.8572
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1426.8572
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1426.8572
$ws.Range("N107").Value = -5266.8572
$ws.Range("H110").Value = 44954.25
$ws.Range("J110").Value = 44954.25
$ws.Range("L110").Value = 44954.25
$ws.Range("N110").Value = -53134.25
$ws.Range("H117").Value = 123003
$ws.Range("J117").Value = 123003
$ws.Range("L117").Value = 123003
$ws.Range("N117").Value = -132181
$ws.Range("H119").Value = 50000
$ws.Range("I119").Value = 50000
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 50000
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -45162
$ws.Range("H120").Value = 79845.5
$ws.Range("J120").Value = 79000
$ws.Range("L120").Value = 79000
$ws.Range("N120").Value = -88676
$ws.Range("H121").Value = 69989
$ws.Range("J121").Value = 69989
$ws.Range("L121").Value = 69989
$ws.Range("N121").Value = -73483
$ws.Range("H141").Value = 98483.164
$ws.Range("J141").Value = 109068.8
$ws.Range("L141").Value = 109068.8
$ws.Range("N141").Value = -119428.8
$ws.Range("M107").ClearContents()
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 931.6667
$ws.Range("I16").Value = 897
$ws.Range("K16").Value = 897
$ws.Range("M16").Value = -610
$ws.Range("H31").Value = 5137.241
$ws.Range("I31").Value = 2746.875
$ws.Range("J31").Value = 6047.857
$ws.Range("K31").Value = 2746.875
$ws.Range("L31").Value = 6047.857
$ws.Range("M31").Value = -2451.875
$ws.Range("N31").Value = -6637.857
$ws.Range("H34").Value = 5137.241
$ws.Range("I34").Value = 2746.875
$ws.Range("J34").Value = 6047.857
$ws.Range("K34").Value = 2746.875
$ws.Range("L34").Value = 6047.857
$ws.Range("M34").Value = -2544.875
$ws.Range("N34").Value = -6451.857
$ws.Range("H58").Value = 2779.2246
$ws.Range("J58").Value = 3760.5
$ws.Range("L58").Value = 3760.5
$ws.Range("N58").Value = -4166.5
$ws.Range("H107").Value = 1591.7222
$ws.Range("I107").Value = 743.8333
$ws.Range("J107").Value = 2015.6666
$ws.Range("K107").Value = 743.8333
$ws.Range("L107").Value = 2015.6666
$ws.Range("M107").Value = 1176.1667
$ws.Range("N107").Value = -5855.6666
$ws.Range("H113").Value = 931.6667
$ws.Range("I113").Value = 897
$ws.Range("K113").Value = 897
$ws.Range("M113").Value = 1273
$ws.Range("H136").Value = 2779.2246
$ws.Range("J136").Value = 3760.5
$ws.Range("L136").Value = 11281.5
$ws.Range("N136").Value = -16381.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 494
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 494
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2964
$ws.Range("N2").Value = -3190
$ws.Range("H18").Value = 4605.7
$ws.Range("I18").Value = 3676.3333
$ws.Range("K18").Value = 11028.9999
$ws.Range("M18").Value = -10859.9999
$ws.Range("H38").Value = 74.5
$ws.Range("J38").Value = 87
$ws.Range("L38").Value = 261
$ws.Range("N38").Value = -955
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28609.715
$ws.Range("I2").Value = 53.6
$ws.Range("J2").Value = 100000
$ws.Range("K2").Value = 53.6
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = 59.4
$ws.Range("N2").Value = -100226
$ws.Range("H14").Value = 1250000
$ws.Range("I14").Value = 1250000
$ws.Range("K14").Value = 1250000
$ws.Range("M14").Value = -1249832
$ws.Range("H107").Value = 1711.75
$ws.Range("I107").Value = 1450
$ws.Range("J107").Value = 1973.5
$ws.Range("K107").Value = 1450
$ws.Range("L107").Value = 1973.5
$ws.Range("M107").Value = 470
$ws.Range("N107").Value = -5813.5
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3555.7144
$ws.Range("I46").Value = 411.75
$ws.Range("J46").Value = 4295.4707
$ws.Range("K46").Value = 411.75
$ws.Range("L46").Value = 4295.4707
$ws.Range("M46").Value = -223.75
$ws.Range("N46").Value = -4671.4707
$ws.Range("H82").Value = 2943.8333
$ws.Range("I82").Value = 2877.5293
$ws.Range("J82").Value = 3104.8572
$ws.Range("K82").Value = 2877.5293
$ws.Range("L82").Value = 3104.8572
$ws.Range("M82").Value = -2516.5293
$ws.Range("N82").Value = -3826.8572
$ws.Range("H85").Value = 2943.8333
$ws.Range("I85").Value = 2877.5293
$ws.Range("J85").Value = 3104.8572
$ws.Range("K85").Value = 2877.5293
$ws.Range("L85").Value = 3104.8572
$ws.Range("M85").Value = -1629.5293
$ws.Range("N85").Value = -5600.8572
$ws.Range("H101").Value = 50772.43
$ws.Range("J101").Value = 50772.43
$ws.Range("L101").Value = 50772.43
$ws.Range("N101").Value = -57262.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 76923464
$ws.Range("I14").Value = 125000000
$ws.Range("J14").Value = 1005
$ws.Range("K14").Value = 125000000
$ws.Range("L14").Value = 1005
$ws.Range("M14").Value = -124999832
$ws.Range("N14").Value = -1341
$ws.Range("H123").Value = 112000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 112000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 112000
$ws.Range("N123").Value = -121800
$ws.Range("H131").Value = 166999.5
$ws.Range("J131").Value = 166999.5
$ws.Range("L131").Value = 166999.5
$ws.Range("N131").Value = -177079.5
$ws.Range("H132").Value = 40951.08
$ws.Range("I132").Value = 44063.707
$ws.Range("J132").Value = 3599.5
$ws.Range("K132").Value = 132191.121
$ws.Range("L132").Value = 10798.5
$ws.Range("M132").Value = -129661.121
$ws.Range("N132").Value = -15858.5
